$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "B401[PD]  /  "
$ws.Range("B4").Value = "CS402[AP]  /  "
$ws.Range("C4").Value = "CS403[MDu]  /  "
$ws.Range("D4").Value = "M401[SMa]  /  "
$ws.Range("E4").Value = "CS401[AH]  /  "
$ws.Range("F4").Value = "CH401[SC]  /  "
$ws.Range("G4").Value = "M401(T)[DC, GY]  /  "
$ws.Range("A6").Value = "CH401[SC]  /  "
$ws.Range("B6").Value = "M401[SMa]  /  "
$ws.Range("C6").Value = "Free Period!"
$ws.Range("D6").Value = "Free Period!"
$ws.Range("E6").Value = "Free Period!"
$ws.Range("F6").Value = "CS401[AH]  /  "
$ws.Range("G6").Value = "Free Period!"
$ws.Range("A8").Value = "CS403[MDu]  /  "
$ws.Range("B8").Value = "CS402[AP]  /  "
$ws.Range("C8").Value = "Free Period!"
$ws.Range("D8").Value = "M401[SMa]  /  "
$ws.Range("E8").Value = "Free Period!"
$ws.Range("F8").Value = "Free Period!"
$ws.Range("G8").Value = "Free Period!"
$ws.Range("A10").Value = "B401[PD]  /  "
$ws.Range("B10").Value = "CS402[AP]  /  "
$ws.Range("C10").Value = "CH401[SC]  /  "
$ws.Range("D10").Value = "CS491[AH, MDu]  /  CS492[AP, BDu]"
$ws.Range("E10").Value = "CS491[AH, MDu]  /  CS492[AP, BDu]"
$ws.Range("F10").Value = "CS491[AH, MDu]  /  CS492[AP, BDu]"
$ws.Range("G10").Value = "CS491[AH, MDu]  /  CS492[AP, BDu]"
$ws.Range("A12").Value = "B401[PD]  /  "
$ws.Range("B12").Value = "CS401[AH]  /  "
$ws.Range("C12").Value = "CS403[MDu]  /  "
$ws.Range("D12").Value = "CS492[AP, BDu]  /  CS491[AH, MDu]"
$ws.Range("E12").Value = "CS492[AP, BDu]  /  CS491[AH, MDu]"
$ws.Range("F12").Value = "CS492[AP, BDu]  /  CS491[AH, MDu]"
$ws.Range("G12").Value = "CS492[AP, BDu]  /  CS491[AH, MDu]"
$ws.Range("A16").Value = "Free Period!"
$ws.Range("B16").Value = "Free Period!"
$ws.Range("C16").Value = "Free Period!"
$ws.Range("D16").Value = "Free Period!"
$ws.Range("E16").Value = "CH401[SC]  /  "
$ws.Range("F16").Value = "IT402[SSR]  /  "
$ws.Range("G16").Value = "IT403[SU]  /  "
$ws.Range("A18").Value = "IT491[SSR, KDa]  /  IT492[RCh, RG]"
$ws.Range("B18").Value = "IT491[SSR, KDa]  /  IT492[RCh, RG]"
$ws.Range("C18").Value = "IT491[SSR, KDa]  /  IT492[RCh, RG]"
$ws.Range("D18").Value = "IT491[SSR, KDa]  /  IT492[RCh, RG]"
$ws.Range("E18").Value = "B401[PD]  /  "
$ws.Range("F18").Value = "IT402[SSR]  /  "
$ws.Range("G18").Value = "M401[ARC]  /  "
$ws.Range("A20").Value = "Free Period!"
$ws.Range("B20").Value = "Free Period!"
$ws.Range("C20").Value = "IT401[RCh]  /  "
$ws.Range("D20").Value = "IT492[RCh, RG]  /  IT491[SSR, KDa]"
$ws.Range("E20").Value = "IT492[RCh, RG]  /  IT491[SSR, KDa]"
$ws.Range("F20").Value = "IT492[RCh, RG]  /  IT491[SSR, KDa]"
$ws.Range("G20").Value = "IT492[RCh, RG]  /  IT491[SSR, KDa]"
$ws.Range("A22").Value = "M401[ARC]  /  "
$ws.Range("B22").Value = "CH401[SC]  /  "
$ws.Range("C22").Value = "IT401[RCh]  /  "
$ws.Range("D22").Value = "Free Period!"
$ws.Range("E22").Value = "B401[PD]  /  "
$ws.Range("F22").Value = "IT402[SSR]  /  "
$ws.Range("G22").Value = "IT403[SU]  /  "
$ws.Range("A24").Value = "M401[ARC]  /  "
$ws.Range("B24").Value = "CH401[SC]  /  "
$ws.Range("C24").Value = "IT401[RCh]  /  "
$ws.Range("D24").Value = "Free Period!"
$ws.Range("E24").Value = "Free Period!"
$ws.Range("F24").Value = "B401[PD]  /  "
$ws.Range("G24").Value = "IT403[SU]  /  "
$ws.Range("A28").Value = "ECE403[PP]  /  "
$ws.Range("B28").Value = "HU491[KB]  /  M491[SLa, SRC]"
$ws.Range("C28").Value = "HU491[KB]  /  M491[SLa, SRC]"
$ws.Range("D28").Value = "B401[PD]  /  "
$ws.Range("E28").Value = "ECE401[SMC]  /  "
$ws.Range("F28").Value = "M401[SLa]  /  "
$ws.Range("G28").Value = "ECE402[SG]  /  "
$ws.Range("A30").Value = "ECE404[SDe]  /  "
$ws.Range("B30").Value = "ECE493[PP, BC]  /  ECE491[PC, SMC]"
$ws.Range("C30").Value = "ECE493[PP, BC]  /  ECE491[PC, SMC]"
$ws.Range("D30").Value = "ECE493[PP, BC]  /  ECE491[PC, SMC]"
$ws.Range("E30").Value = "ECE491[PC, SMC]  /  ECE492[DK, JA]"
$ws.Range("F30").Value = "ECE491[PC, SMC]  /  ECE492[DK, JA]"
$ws.Range("G30").Value = "ECE491[PC, SMC]  /  ECE492[DK, JA]"
$ws.Range("A32").Value = "ECE403[PP]  /  "
$ws.Range("B32").Value = "ECE492[DK, JA]  /  ECE493[PP, BC]"
$ws.Range("C32").Value = "ECE492[DK, JA]  /  ECE493[PP, BC]"
$ws.Range("D32").Value = "ECE492[DK, JA]  /  ECE493[PP, BC]"
$ws.Range("E32").Value = "ECE401[SMC]  /  "
$ws.Range("F32").Value = "Free Period!"
$ws.Range("G32").Value = "ECE404[SDe]  /  "
$ws.Range("A34").Value = "ECE403[PP]  /  "
$ws.Range("B34").Value = "B401[PD]  /  "
$ws.Range("C34").Value = "ECE402[SG]  /  "
$ws.Range("D34").Value = "Free Period!"
$ws.Range("E34").Value = "Free Period!"
$ws.Range("F34").Value = "M401[SLa]  /  "
$ws.Range("G34").Value = "ECE404[SDe]  /  "
$ws.Range("A36").Value = "M491[SLa, SRC]  /  HU491[KB]"
$ws.Range("B36").Value = "M491[SLa, SRC]  /  HU491[KB]"
$ws.Range("C36").Value = "M401[SLa]  /  "
$ws.Range("D36").Value = "ECE401[SMC]  /  "
$ws.Range("E36").Value = "B401[PD]  /  "
$ws.Range("F36").Value = "Free Period!"
$ws.Range("G36").Value = "ECE402[SG]  /  "
$ws.Range("A40").Value = "EE492[RND, JA]  /  EE491[KR, IB]"
$ws.Range("B40").Value = "EE492[RND, JA]  /  EE491[KR, IB]"
$ws.Range("C40").Value = "EE492[RND, JA]  /  EE491[KR, IB]"
$ws.Range("D40").Value = "EE404[ArD]  /  "
$ws.Range("E40").Value = "EE494[BDC, NCS]  /  EE493[ArD, SDG]"
$ws.Range("F40").Value = "EE494[BDC, NCS]  /  EE493[ArD, SDG]"
$ws.Range("G40").Value = "EE494[BDC, NCS]  /  EE493[ArD, SDG]"
$ws.Range("A42").Value = "HU401[ACh]  /  "
$ws.Range("B42").Value = "EE403[KR]  /  "
$ws.Range("C42").Value = "EE404[ArD]  /  "
$ws.Range("D42").Value = "EE401[BDC]  /  "
$ws.Range("E42").Value = "EE491[KR, IB]  /  EE494[BDC, NCS]"
$ws.Range("F42").Value = "EE491[KR, IB]  /  EE494[BDC, NCS]"
$ws.Range("G42").Value = "EE491[KR, IB]  /  EE494[BDC, NCS]"
$ws.Range("A44").Value = "CH401[PD]  /  "
$ws.Range("B44").Value = "HU401[ACh]  /  "
$ws.Range("C44").Value = "EE402[RND]  /  "
$ws.Range("D44").Value = "EE401[BDC]  /  "
$ws.Range("E44").Value = "EE493[ArD, SDG]  /  EE492[RND, JA]"
$ws.Range("F44").Value = "EE493[ArD, SDG]  /  EE492[RND, JA]"
$ws.Range("G44").Value = "EE493[ArD, SDG]  /  EE492[RND, JA]"
$ws.Range("A46").Value = "CH401[SC]  /  "
$ws.Range("B46").Value = "Free Period!"
$ws.Range("C46").Value = "EE402[RND]  /  "
$ws.Range("D46").Value = "Free Period!"
$ws.Range("E46").Value = "HU401[ACh]  /  "
$ws.Range("F46").Value = "Free Period!"
$ws.Range("G46").Value = "EE403[KR]  /  "
$ws.Range("A48").Value = "CH401[SC]  /  "
$ws.Range("B48").Value = "EE401[BDC]  /  "
$ws.Range("C48").Value = "EE404[ArD]  /  "
$ws.Range("D48").Value = "Free Period!"
$ws.Range("E48").Value = "Free Period!"
$ws.Range("F48").Value = "EE402[RND]  /  "
$ws.Range("G48").Value = "EE403[KR]  /  "
